$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells D1/E1, copying the header style from an existing
# header cell (C1) so they get the same bold/centered/bordered formatting
# that the other headers use.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "price_predicted"

$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "residuals"

# Fill in the predicted price and residual values for each data row.
$ws.Range("D2").Value = 41.11111111111111
$ws.Range("E2").Value = 6.111111111111114

$ws.Range("D3").Value = 46.38888888888889
$ws.Range("E3").Value = -13.61111111111111

$ws.Range("D4").Value = 49.02777777777779
$ws.Range("E4").Value = 29.02777777777779

$ws.Range("D5").Value = 43.75
$ws.Range("E5").Value = -6.25

$ws.Range("D6").Value = 46.38888888888889
$ws.Range("E6").Value = -3.611111111111107

$ws.Range("D7").Value = 49.02777777777779
$ws.Range("E7").Value = -5.972222222222214

$ws.Range("D8").Value = 54.30555555555556
$ws.Range("E8").Value = -5.694444444444443
